# Moderate overlay ready for review:
#  - rename the worksheet from the "high" overlay name to the "moderate" one
#  - keep the hidden _FilterDatabase defined name in sync with the new sheet name
#  - move the active selection (in the frozen bottom-right pane) to D4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (also updates the sheet reference used by the
# _xlnm._FilterDatabase defined name automatically).
$ws.Name = "cmsars3.1moderate_w2019overlay"

# Re-assert the hidden filter-database defined name so the sheet name in its
# formula is quoted exactly like the rest of the workbook.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='cmsars3.1moderate_w2019overlay'!`$A`$1:`$U`$304"
    }
}

# Move the selection to D4 (falls in the frozen bottom-right pane).
$ws.Range("D4").Select()
